$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its title cell / sharedString text
$ws.Name = "Through 2022-03-04"
$ws.Range("A4").Value = "March (through 03-04)"

# Update February (row 3) value for 2022 (column I)
$ws.Range("I3").Value = 142

# Update March (row 4) values across all year columns
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 20

# Update Total (row 5) values across all year columns
$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 91
$ws.Range("D5").Value = 135
$ws.Range("E5").Value = 145
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 148
$ws.Range("H5").Value = 354
$ws.Range("I5").Value = 321
